$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85; this shifts existing rows 85..174 down to 86..175
$ws.Rows("85:85").Insert()

# Populate the newly inserted row 85 with the new data record
$ws.Range("A85").Value = 5
$ws.Range("B85").Value = "Macroferia Regional de Talca"
$ws.Range("C85").Value = "Maule"
$ws.Range("D85").Value = 44589
$ws.Range("E85").Value = 7
$ws.Range("F85").Value = 100112021
$ws.Range("G85").Value = "Ají"
$ws.Range("H85").Value = "Americana (o)"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 150
$ws.Range("K85").Value = 18000
$ws.Range("L85").Value = 18000
$ws.Range("M85").Value = 18000
$ws.Range("N85").Value = "`$/saco 25 kilos"
$ws.Range("O85").Value = "Región del Maule"
$ws.Range("P85").Value = 720
$ws.Range("Q85").Value = 25
$ws.Range("R85").Value = "Hortaliza"
